# Update the "想去人数" (F column) figures across the 展览, 演出 and 全部类型
# sheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsShow       = $wb.Worksheets.Item("演出")
$wsAll        = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet (sheet1) ---
$exhibitionUpdates = @{
    2  = 791
    3  = 528
    4  = 276
    5  = 480
    6  = 1127
    8  = 25
    9  = 111
    10 = 109
    11 = 1127
    14 = 770
    15 = 801
    17 = 40
    18 = 62
    19 = 665
    20 = 169
    21 = 1711
    22 = 2223
    23 = 608
    25 = 1851
    26 = 302
    27 = 2700
    28 = 487
    29 = 78
    30 = 673
    31 = 129
    32 = 95
    33 = 94
    34 = 939
    35 = 1660
    36 = 306
    38 = 532
    39 = 146
    40 = 113
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# --- 演出 sheet (sheet2) ---
$showUpdates = @{
    3 = 134
    8 = 7
}

foreach ($row in $showUpdates.Keys) {
    $wsShow.Range("F$row").Value = $showUpdates[$row]
}

# --- 全部类型 sheet (sheet4) ---
$allUpdates = @{
    3  = 791
    4  = 528
    5  = 276
    6  = 480
    7  = 1127
    9  = 25
    10 = 111
    11 = 109
    12 = 1127
    14 = 770
    15 = 801
    17 = 134
    18 = 134
    20 = 40
    22 = 62
    23 = 169
    24 = 1711
    25 = 2223
    26 = 608
    29 = 1851
    31 = 2700
    32 = 487
    33 = 7
    36 = 78
    38 = 673
    39 = 129
    40 = 95
    41 = 94
    42 = 939
    43 = 1660
    45 = 306
    46 = 532
    47 = 146
    48 = 113
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
